$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 17.040744
$ws.Range("D3").Value = 69.550248
